$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row total correct count
$ws.Range("B11").Value = 5

# Update "Total" row correct count and corr/total marks display
$ws.Range("B12").Value = 115
$ws.Range("E12").Value = "115/140"
